$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2606.9285
$ws.Range("J2").Value = 5195
$ws.Range("L2").Value = 5195
$ws.Range("N2").Value = -5421
$ws.Range("H40").Value = 2231.5557
$ws.Range("I40").Value = 2105.3845
$ws.Range("K40").Value = 2105.3845
$ws.Range("M40").Value = -1930.3845
$ws.Range("H62").Value = 6863.5713
$ws.Range("I62").Value = 6409
$ws.Range("K62").Value = 6409
$ws.Range("M62").Value = -5785
$ws.Range("H65").Value = 6863.5713
$ws.Range("I65").Value = 6409
$ws.Range("K65").Value = 32045
$ws.Range("M65").Value = -28925
$ws.Range("H76").Value = 8332.666999999999
$ws.Range("J76").Value = 7499
$ws.Range("L76").Value = 7499
$ws.Range("N76").Value = -8129
$ws.Range("H79").Value = 8332.666999999999
$ws.Range("J79").Value = 7499
$ws.Range("L79").Value = 7499
$ws.Range("N79").Value = -9683
$ws.Range("H86").Value = 3077.7144
$ws.Range("I86").Value = 1149.25
$ws.Range("K86").Value = 1149.25
$ws.Range("M86").Value = -26.25
$ws.Range("H89").Value = 3077.7144
$ws.Range("I89").Value = 1149.25
$ws.Range("K89").Value = 5746.25
$ws.Range("M89").Value = -130.25
$ws.Range("H112").Value = 1530.8667
$ws.Range("J112").Value = 1531.931
$ws.Range("L112").Value = 4595.793
$ws.Range("N112").Value = -6811.793
$ws.Range("H113").Value = 4093.5
$ws.Range("J113").Value = 3958.3333
$ws.Range("L113").Value = 3958.3333
$ws.Range("N113").Value = -10466.3333
$ws.Range("H138").Value = 5778.34
$ws.Range("J138").Value = 6187.564
$ws.Range("L138").Value = 18562.692
$ws.Range("N138").Value = -28842.692
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8562.966
$ws.Range("I32").Value = 6619.0356
$ws.Range("K32").Value = 6619.0356
$ws.Range("M32").Value = -6332.0356
$ws.Range("H61").Value = 7943.1113
$ws.Range("J61").Value = 7998.5
$ws.Range("L61").Value = 7998.5
$ws.Range("N61").Value = -8422.5
$ws.Range("H74").Value = 2917.125
$ws.Range("I74").Value = 2875.077
$ws.Range("J74").Value = 3099.3333
$ws.Range("K74").Value = 2875.077
$ws.Range("L74").Value = 3099.3333
$ws.Range("M74").Value = -2001.077
$ws.Range("N74").Value = -4847.3333
$ws.Range("H77").Value = 2917.125
$ws.Range("I77").Value = 2875.077
$ws.Range("J77").Value = 3099.3333
$ws.Range("K77").Value = 14375.385
$ws.Range("L77").Value = 15496.6665
$ws.Range("M77").Value = -10007.385
$ws.Range("N77").Value = -24232.6665
$ws.Range("H102").Value = 2304.0833
$ws.Range("I102").Value = 2365
$ws.Range("K102").Value = 2365
$ws.Range("M102").Value = -743
$ws.Range("H122").Value = 2602.6924
$ws.Range("I122").Value = 2333.1
$ws.Range("J122").Value = 3501.3333
$ws.Range("K122").Value = 6999.299999999999
$ws.Range("L122").Value = 10503.9999
$ws.Range("M122").Value = -4549.299999999999
$ws.Range("N122").Value = -15403.9999
$ws.Range("H132").Value = 7777.222
$ws.Range("I132").Value = 3665.8333
$ws.Range("K132").Value = 10997.4999
$ws.Range("M132").Value = -8467.499899999999
$ws.Range("H136").Value = 7943.1113
$ws.Range("J136").Value = 7998.5
$ws.Range("L136").Value = 23995.5
$ws.Range("N136").Value = -29095.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1307
$ws.Range("I22").Value = 1307
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1307
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1134
$ws.Range("N22").ClearContents()
$ws.Range("H86").Value = 2812.0833
$ws.Range("I86").Value = 2502.4
$ws.Range("J86").Value = 4360.5
$ws.Range("K86").Value = 2502.4
$ws.Range("L86").Value = 4360.5
$ws.Range("M86").Value = -1379.4
$ws.Range("N86").Value = -6606.5
$ws.Range("H89").Value = 2812.0833
$ws.Range("I89").Value = 2502.4
$ws.Range("J89").Value = 4360.5
$ws.Range("K89").Value = 12512
$ws.Range("L89").Value = 21802.5
$ws.Range("M89").Value = -6896
$ws.Range("N89").Value = -33034.5
$ws.Range("H99").Value = 3324.1
$ws.Range("I99").Value = 2775.7144
$ws.Range("J99").Value = 4603.6665
$ws.Range("K99").Value = 2775.7144
$ws.Range("L99").Value = 4603.6665
$ws.Range("M99").Value = -1277.7144
$ws.Range("N99").Value = -7599.6665
$ws.Range("H134").Value = 3845.2727
$ws.Range("I134").Value = 3212.1428
$ws.Range("J134").Value = 4953.25
$ws.Range("K134").Value = 9636.428400000001
$ws.Range("L134").Value = 14859.75
$ws.Range("M134").Value = -7101.428400000001
$ws.Range("N134").Value = -19929.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H37").Value = 7999
$ws.Range("I37").Value = 5999
$ws.Range("J37").Value = 9999
$ws.Range("K37").Value = 5999
$ws.Range("L37").Value = 9999
$ws.Range("M37").Value = -5892
$ws.Range("N37").Value = -10213
$ws.Range("H58").Value = 1399
$ws.Range("I58").Value = 1399
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1399
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1196
$ws.Range("N58").ClearContents()
$ws.Range("H86").Value = 8608.143
$ws.Range("J86").Value = 6179
$ws.Range("L86").Value = 6179
$ws.Range("N86").Value = -8425
$ws.Range("H89").Value = 8608.143
$ws.Range("J89").Value = 6179
$ws.Range("L89").Value = 30895
$ws.Range("N89").Value = -42127
$ws.Range("H103").Value = 10424.333
$ws.Range("I103").Value = 7709.4
$ws.Range("J103").Value = 23999
$ws.Range("K103").Value = 7709.4
$ws.Range("L103").Value = 23999
$ws.Range("M103").Value = -6537.4
$ws.Range("N103").Value = -26343
$ws.Range("H132").Value = 1487.2
$ws.Range("I132").Value = 970.25
$ws.Range("J132").Value = 3555
$ws.Range("K132").Value = 2910.75
$ws.Range("L132").Value = 10665
$ws.Range("M132").Value = -380.75
$ws.Range("N132").Value = -15725
$ws.Range("H134").Value = 2653.5715
$ws.Range("I134").Value = 2653.5715
$ws.Range("K134").Value = 7960.7145
$ws.Range("M134").Value = -5425.7145
$ws.Range("H136").Value = 1399
$ws.Range("I136").Value = 1399
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4197
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1647
$ws.Range("N136").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 130.77777
$ws.Range("J12").Value = 214.8
$ws.Range("L12").Value = 644.4000000000001
$ws.Range("N12").Value = -990.4000000000001
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H63").Value = 17721.25
$ws.Range("I63").Value = 17328.334
$ws.Range("J63").Value = 18900
$ws.Range("K63").Value = 51985.00199999999
$ws.Range("L63").Value = 56700
$ws.Range("M63").Value = -51236.00199999999
$ws.Range("N63").Value = -58198
$ws.Range("H66").Value = 17721.25
$ws.Range("I66").Value = 17328.334
$ws.Range("J66").Value = 18900
$ws.Range("K66").Value = 155955.006
$ws.Range("L66").Value = 170100
$ws.Range("M66").Value = -152211.006
$ws.Range("N66").Value = -177588
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 376
$ws.Range("I107").Value = 390.33334
$ws.Range("J107").Value = 247
$ws.Range("K107").Value = 390.33334
$ws.Range("L107").Value = 247
$ws.Range("M107").Value = 1529.66666
$ws.Range("N107").Value = -4087
$ws.Range("H122").Value = 4524
$ws.Range("I122").Value = 4524
$ws.Range("K122").Value = 13572
$ws.Range("M122").Value = -11122
$ws.Range("H126").Value = 4003.6667
$ws.Range("J126").Value = 4337.3335
$ws.Range("L126").Value = 13012.0005
$ws.Range("N126").Value = -17952.0005
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1850.9048
$ws.Range("I82").Value = 1465.8334
$ws.Range("J82").Value = 2364.3333
$ws.Range("K82").Value = 1465.8334
$ws.Range("L82").Value = 2364.3333
$ws.Range("M82").Value = -1104.8334
$ws.Range("N82").Value = -3086.3333
$ws.Range("H85").Value = 1850.9048
$ws.Range("I85").Value = 1465.8334
$ws.Range("J85").Value = 2364.3333
$ws.Range("K85").Value = 1465.8334
$ws.Range("L85").Value = 2364.3333
$ws.Range("M85").Value = -217.8334
$ws.Range("N85").Value = -4860.3333
$ws.Range("H122").Value = 7598.3213
$ws.Range("I122").Value = 8393.8125
$ws.Range("J122").Value = 6537.6665
$ws.Range("K122").Value = 25181.4375
$ws.Range("L122").Value = 19612.9995
$ws.Range("M122").Value = -22731.4375
$ws.Range("N122").Value = -24512.9995
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 54500
$ws.Range("J92").Value = 54500
$ws.Range("L92").Value = 54500
$ws.Range("N92").Value = -59492
